# Re-verify Feature #138 (Sortie search and filtering) asset-report refresh:
#  - bump the "Generated" timestamp in A4
#  - the CRIIS-010 / Navigation Unit asset (row 17) is gone, so the total
#    asset count in A5 drops from 10 to 9
#  - the row 17 record itself is removed; everything below (the CUI footer,
#    previously on row 19) shifts up to row 18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Generated" timestamp (A4)
$ws.Range("A4").Value = "Generated: 2026-01-20 11:39:01Z"

# Update the total asset count (A5): 10 -> 9
$ws.Range("A5").Value = "Total Assets: 9"

# Remove the CRIIS-010 (Navigation Unit) data row entirely. Deleting the
# whole row shifts every row below it up by one -- the CUI footer row
# (formerly row 19) becomes row 18, the sheet dimension shrinks to
# A1:L18, and the footer's merged cell range follows it to A18:L18.
$ws.Rows(17).Delete()
